$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "195.80" or "1.002"
# are not auto-converted to numbers (matches original inline-string text cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.984.11"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.640.43"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.82%  "

$ws.Range("D5").Value = "214.87"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").Value = "0.5067"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "0.06355"
$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("D10").Value = "19.83"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").Value = "0.07753"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").Value = "4.284"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "1.642.44"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").Value = "0.5485"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "64.32"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₅7731"
$ws.Range("E16").Value = "  -2.18%  "

$ws.Range("D17").Value = "26.014.51"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "195.80"
$ws.Range("E19").Value = "  -1.51%  "

$ws.Range("D20").Value = "4.440"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").Value = "9.968"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("D22").Value = "6.120"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  +1.14%  "

$ws.Range("D25").Value = "143.25"
$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").Value = "0.1256"
$ws.Range("E26").Value = "  +9.27%  "

$ws.Range("D27").Value = "6.875"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").Value = "15.60"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").Value = "0.04891"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("D31").Value = "3.272"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").Value = "3.210"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").Value = "1.547"

$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("D35").Value = "0.9159"
$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").Value = "0.5547"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").Value = "1.096.13"
$ws.Range("E38").Value = "  -3.51%  "

$ws.Range("D39").Value = "0.01571"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("D41").Value = "5.612"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").Value = "0.8043"
$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("D43").Value = "98.97"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("D44").Value = "0.0₈121"
$ws.Range("E44").Value = "  -3.78%  "

$ws.Range("D45").Value = "1.781.49"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "0.4535"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "55.34"
$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("D48").Value = "1.002"

$ws.Range("D49").Value = "0.05198"
$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("D50").Value = "7.555"
$ws.Range("E50").Value = "  +1.78%  "

$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -0.22%  "
